$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (recomputed after adding more games
# and speeding up the simulate-game logic).
$ws.Range("B2").Value = 0.2078977932636469
$ws.Range("C2").Value = 0.537746806039489
$ws.Range("J2").Value = 0.0116144018583043
$ws.Range("P2").Value = 0.1556329849012776
$ws.Range("S2").Value = 0.08710801393728224
$ws.Range("B3").Value = 0.02691511387163561
$ws.Range("C3").Value = 0.03933747412008282
$ws.Range("J3").Value = 0.03105590062111801
$ws.Range("P3").Value = 0.6749482401656315
$ws.Range("S3").Value = 0.2277432712215321
$ws.Range("J4").Value = 0.02654867256637168
$ws.Range("P4").Value = 0.6637168141592921
$ws.Range("S4").Value = 0.3097345132743363
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.1014975041597338
$ws.Range("D6").Value = 0.01663893510815308
$ws.Range("F6").Value = 0.07986688851913477
$ws.Range("J6").Value = 0.2412645590682196
$ws.Range("O6").Value = 0.01830282861896839
$ws.Range("Q6").Value = 0.1397670549084858
$ws.Range("R6").Value = 0.04991680532445923
$ws.Range("S6").Value = 0.3527454242928453
$ws.Range("B7").Value = 0.118546845124283
$ws.Range("D7").Value = 0.02103250478011472
$ws.Range("E7").Value = 0.001912045889101338
$ws.Range("F7").Value = 0.06500956022944551
$ws.Range("J7").Value = 0.1395793499043977
$ws.Range("O7").Value = 0.0248565965583174
$ws.Range("Q7").Value = 0.1472275334608031
$ws.Range("R7").Value = 0.07074569789674952
$ws.Range("S7").Value = 0.4110898661567878
$ws.Range("B8").Value = 0.09272581934452438
$ws.Range("D8").Value = 0.01438848920863309
$ws.Range("E8").Value = 0.001598721023181455
$ws.Range("F8").Value = 0.05995203836930456
$ws.Range("J8").Value = 0.1239008792965628
$ws.Range("O8").Value = 0.01278976818545164
$ws.Range("Q8").Value = 0.17585931254996
$ws.Range("R8").Value = 0.09592326139088729
$ws.Range("S8").Value = 0.4228617106314948
$ws.Range("B9").Value = 0.0945945945945946
$ws.Range("D9").Value = 0.02316602316602316
$ws.Range("F9").Value = 0.06563706563706563
$ws.Range("J9").Value = 0.1177606177606178
$ws.Range("O9").Value = 0.01544401544401544
$ws.Range("Q9").Value = 0.1814671814671815
$ws.Range("R9").Value = 0.09073359073359073
$ws.Range("S9").Value = 0.4111969111969112
$ws.Range("B10").Value = 0.1116781747223056
$ws.Range("D10").Value = 0.01981386970879616
$ws.Range("E10").Value = 0.00150105073551486
$ws.Range("F10").Value = 0.07054938456919845
$ws.Range("J10").Value = 0.1371960372260582
$ws.Range("O10").Value = 0.01471029720804563
$ws.Range("Q10").Value = 0.1873311317922546
$ws.Range("R10").Value = 0.07325127589312519
$ws.Range("S10").Value = 0.3839687781447013
$ws.Range("G11").Value = 0.1537484116899619
$ws.Range("J11").Value = 0.07623888182973317
$ws.Range("K11").Value = 0.1944091486658196
$ws.Range("L11").Value = 0.5527318932655655
$ws.Range("S11").Value = 0.02287166454891995
$ws.Range("G12").Value = 0.7347826086956522
$ws.Range("J12").Value = 0.1717391304347826
$ws.Range("K12").Value = 0.008695652173913044
$ws.Range("L12").Value = 0.04130434782608695
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("G13").Value = 0.7352941176470589
$ws.Range("J13").Value = 0.2058823529411765
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.02843601895734597
$ws.Range("H15").Value = 0.1658767772511848
$ws.Range("I15").Value = 0.08688783570300158
$ws.Range("J15").Value = 0.368088467614534
$ws.Range("K15").Value = 0.06161137440758294
$ws.Range("M15").Value = 0.009478672985781991
$ws.Range("O15").Value = 0.08372827804107424
$ws.Range("S15").Value = 0.1958925750394945
$ws.Range("F16").Value = 0.01724137931034483
$ws.Range("H16").Value = 0.1934865900383142
$ws.Range("I16").Value = 0.06321839080459771
$ws.Range("J16").Value = 0.3946360153256705
$ws.Range("K16").Value = 0.1053639846743295
$ws.Range("M16").Value = 0.01915708812260536
$ws.Range("O16").Value = 0.08045977011494253
$ws.Range("S16").Value = 0.1264367816091954
$ws.Range("F17").Value = 0.01931922723091076
$ws.Range("H17").Value = 0.1895124195032199
$ws.Range("I17").Value = 0.08463661453541858
$ws.Range("J17").Value = 0.3937442502299908
$ws.Range("K17").Value = 0.1021159153633855
$ws.Range("M17").Value = 0.01655933762649494
$ws.Range("N17").Value = 0.001839926402943882
$ws.Range("O17").Value = 0.07543698252069918
$ws.Range("S17").Value = 0.1168353265869365
$ws.Range("F18").Value = 0.02489626556016597
$ws.Range("H18").Value = 0.1659751037344398
$ws.Range("I18").Value = 0.1120331950207469
$ws.Range("J18").Value = 0.4211618257261411
$ws.Range("K18").Value = 0.09336099585062241
$ws.Range("M18").Value = 0.01867219917012448
$ws.Range("O18").Value = 0.07883817427385892
$ws.Range("S18").Value = 0.08506224066390042
$ws.Range("F19").Value = 0.01755926251097454
$ws.Range("H19").Value = 0.2197834357623646
$ws.Range("I19").Value = 0.08282118817676325
$ws.Range("J19").Value = 0.3652326602282704
$ws.Range("K19").Value = 0.1082821188176763
$ws.Range("M19").Value = 0.01785191688615745
$ws.Range("N19").Value = 0.001755926251097454
$ws.Range("O19").Value = 0.07052970441908106
$ws.Range("S19").Value = 0.1161837869476149
